$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "StatQuery" column (C) for CasesTab/SamplesTab/FilesTab rows was
# rewritten with a corrected Cypher query (commit: "Fixed ICDC breed all
# testcases"). Replace the old long query text in C2:C4 with the new one.
$newStatQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.breed IN ['Chinese Shar-Pei']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

$ws.Range("C2").Value = $newStatQuery
$ws.Range("C3").Value = $newStatQuery
$ws.Range("C4").Value = $newStatQuery

# View state was also updated when the file was last saved (zoomed in,
# scrolled down, selection moved to B4).
$excel.ActiveWindow.Zoom = 85
$ws.Range("B4").Select()
